$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value2 = 6242
$ws.Range("C21").Value2 = 988
$ws.Range("D21").Value2 = 5623695
$ws.Range("E21").Value2 = 900.9444088433195
$ws.Range("F21").Value2 = 8.34924492275646
$ws.Range("G21").Value2 = 4.219409282700415
$ws.Range("H21").Value2 = 28.36637622883689
